$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = 3.075165666666667
$ws.Cells.Item(2, 8).Value = 9.225497000000001
$ws.Cells.Item(2, 9).Value = 0.02641273658732285
$ws.Cells.Item(2, 10).Value = 0.02641273658732285
$ws.Cells.Item(2, 13).Value = 0.2901893333333334
$ws.Cells.Item(2, 14).Value = 0.870568
$ws.Cells.Item(2, 15).Value = 0.03429389578125064
$ws.Cells.Item(2, 16).Value = 0.03429389578125064
$ws.Cells.Item(2, 17).Value = 0.8923802746995557
$ws.Cells.Item(2, 18).Value = 8.031422472296001
$ws.Cells.Item(2, 19).Value = 0.0009057956358232754
$ws.Cells.Item(2, 20).Value = 0.0009057956358232754
$ws.Cells.Item(3, 7).Value = 3.075165666666667
$ws.Cells.Item(3, 8).Value = 9.225497000000001
$ws.Cells.Item(3, 9).Value = 0.02641273658732285
$ws.Cells.Item(3, 10).Value = 0.02641273658732285
$ws.Cells.Item(3, 15).Value = 0.8402845891331153
$ws.Cells.Item(3, 16).Value = 0.8402845891331153
$ws.Cells.Item(3, 17).Value = 21.86550624809378
$ws.Cells.Item(3, 18).Value = 196.789556232844
$ws.Cells.Item(3, 19).Value = 0.02219421551115978
$ws.Cells.Item(3, 20).Value = 0.02219421551115978
$ws.Cells.Item(4, 7).Value = 3.075165666666667
$ws.Cells.Item(4, 8).Value = 9.225497000000001
$ws.Cells.Item(4, 9).Value = 0.02641273658732285
$ws.Cells.Item(4, 10).Value = 0.02641273658732285
$ws.Cells.Item(4, 15).Value = 0.1254215150856341
$ws.Cells.Item(4, 16).Value = 0.1254215150856341
$ws.Cells.Item(4, 17).Value = 3.263662046425889
$ws.Cells.Item(4, 18).Value = 29.372958417833
$ws.Cells.Item(4, 19).Value = 0.003312725440339793
$ws.Cells.Item(4, 20).Value = 0.003312725440339792
$ws.Cells.Item(5, 9).Value = 0.549422396165273
$ws.Cells.Item(5, 10).Value = 0.5494223961652731
$ws.Cells.Item(5, 13).Value = 0.2901893333333334
$ws.Cells.Item(5, 14).Value = 0.870568
$ws.Cells.Item(5, 15).Value = 0.03429389578125064
$ws.Cells.Item(5, 16).Value = 0.03429389578125064
$ws.Cells.Item(5, 17).Value = 18.562775848504
$ws.Cells.Item(5, 18).Value = 167.064982636536
$ws.Cells.Item(5, 19).Value = 0.01884183439397687
$ws.Cells.Item(5, 20).Value = 0.01884183439397688
$ws.Cells.Item(6, 9).Value = 0.549422396165273
$ws.Cells.Item(6, 10).Value = 0.5494223961652731
$ws.Cells.Item(6, 15).Value = 0.8402845891331153
$ws.Cells.Item(6, 16).Value = 0.8402845891331153
$ws.Cells.Item(6, 19).Value = 0.4616711724222681
$ws.Cells.Item(6, 20).Value = 0.4616711724222682
$ws.Cells.Item(7, 9).Value = 0.549422396165273
$ws.Cells.Item(7, 10).Value = 0.5494223961652731
$ws.Cells.Item(7, 15).Value = 0.1254215150856341
$ws.Cells.Item(7, 16).Value = 0.1254215150856341
$ws.Cells.Item(7, 17).Value = 67.888801143067
$ws.Cells.Item(7, 19).Value = 0.06890938934902803
$ws.Cells.Item(7, 20).Value = 0.06890938934902803
$ws.Cells.Item(8, 9).Value = 0.424164867247404
$ws.Cells.Item(8, 10).Value = 0.4241648672474041
$ws.Cells.Item(8, 13).Value = 0.2901893333333334
$ws.Cells.Item(8, 14).Value = 0.870568
$ws.Cells.Item(8, 15).Value = 0.03429389578125064
$ws.Cells.Item(8, 16).Value = 0.03429389578125064
$ws.Cells.Item(8, 17).Value = 14.33082708036445
$ws.Cells.Item(8, 18).Value = 128.97744372328
$ws.Cells.Item(8, 19).Value = 0.01454626575145049
$ws.Cells.Item(8, 20).Value = 0.01454626575145049
$ws.Cells.Item(9, 9).Value = 0.424164867247404
$ws.Cells.Item(9, 10).Value = 0.4241648672474041
$ws.Cells.Item(9, 15).Value = 0.8402845891331153
$ws.Cells.Item(9, 16).Value = 0.8402845891331153
$ws.Cells.Item(9, 19).Value = 0.3564192011996873
$ws.Cells.Item(9, 20).Value = 0.3564192011996873
$ws.Cells.Item(10, 9).Value = 0.424164867247404
$ws.Cells.Item(10, 10).Value = 0.4241648672474041
$ws.Cells.Item(10, 15).Value = 0.1254215150856341
$ws.Cells.Item(10, 16).Value = 0.1254215150856341
$ws.Cells.Item(10, 17).Value = 52.41148618152111
$ws.Cells.Item(10, 19).Value = 0.05319940029626628
$ws.Cells.Item(10, 20).Value = 0.05319940029626628
